$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 20, shifting rows 20:43 down to 21:44
$ws.Rows("20:20").Insert()

# Populate the newly inserted row 20 with its values
$ws.Range("A20").Value = 10
$ws.Range("B20").Value = "Vega Modelo de Temuco"
$ws.Range("C20").Value = "La Araucanía"
$ws.Range("D20").Value = 44484
$ws.Range("D20").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E20").Value = 9
$ws.Range("F20").Value = 300000001
$ws.Range("G20").Value = "Rabanito"
$ws.Range("H20").Value = "Sin especificar"
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 30
$ws.Range("K20").Value = 7000
$ws.Range("L20").Value = 8000
$ws.Range("M20").Value = 7333
$ws.Range("N20").Value = "$/docena de paquetes"
$ws.Range("O20").Value = "Provincia de Cautín"
$ws.Range("P20").Value = 611
$ws.Range("Q20").Value = 12
$ws.Range("R20").Value = "Hortaliza"
